$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A95").Value = "Create Country"
$ws.Range("B95").Value = "PASSED"
$ws.Range("C95").Value = "chrome"
$ws.Range("D95").Value = "'07.04.23"
